# Rename the "Severidade da NC" header (column C, row 2) to "Criticidade da NC"
# on every worksheet of the workbook, then leave the view state matching the
# author's editing session (active sheet = "Teste", with each sheet's last
# selection reflecting where the edit/navigation left off).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Processo
$ws1.Activate()
$ws1.Range("C2").Value = "Criticidade da NC"
$ws1.Range("C2").Select()

$ws2 = $wb.Worksheets.Item(2)   # Análise
$ws2.Activate()
$ws2.Range("C2").Value = "Criticidade da NC"
$ws2.Range("C2").Select()

$ws3 = $wb.Worksheets.Item(3)   # Projeto
$ws3.Activate()
$ws3.Range("C2").Value = "Criticidade da NC"
$ws3.Range("C5").Select()

$ws4 = $wb.Worksheets.Item(4)   # Codificação
$ws4.Activate()
$ws4.Range("C2").Value = "Criticidade da NC"
$ws4.Range("C2").Select()

$ws5 = $wb.Worksheets.Item(5)   # Teste
$ws5.Activate()
$ws5.Range("C2").Value = "Criticidade da NC"
$ws5.Range("C2").Select()
